$d = $word.ActiveDocument

# Replaces the single paragraph whose current text equals $OldText (including
# its trailing CR) with the literal paragraph markup in $ParaXml, preserving
# exact run structure (e.g. leading empty <w:r/> runs that plain Find/Replace
# or Range.Text assignment would otherwise silently merge away).
function Replace-ParagraphXml($OldText, $ParaXml) {
    $paras = $d.Paragraphs
    $target = $null
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -eq $OldText) {
            $target = $p
            break
        }
    }
    if ($target -eq $null) { return $false }

    $countBefore = $d.Paragraphs.Count

    $rng = $target.Range
    $rng.End = $rng.End - 1   # exclude the paragraph mark itself
    $insertStart = $rng.Start
    $rng.Delete()

    $pkg = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $ParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)

    $countAfter = $d.Paragraphs.Count
    if ($countAfter -gt $countBefore) {
        # InsertXML always mints a fresh paragraph mark for the <w:p> we hand
        # it. That is a no-op everywhere except when replacing the document's
        # very last paragraph, where it leaves one extra empty paragraph
        # behind. Merge it back out by deleting the duplicate mark right
        # after our freshly-inserted paragraph's content.
        for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
            $pp = $d.Paragraphs.Item($i)
            if ($pp.Range.Start -eq $insertStart) {
                $markRng = $d.Range($pp.Range.End - 1, $pp.Range.End)
                $markRng.Delete()
                break
            }
        }
    }

    return $true
}

# 1. Main H1 title.
Replace-ParagraphXml "Play Book of Souls II: El Dorado for Free - Review`r" `
    '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Book of Souls II: El Dorado Free | Adventure-Themed Slot</w:t></w:r></w:p>'

# 2. "What we like" bullet list (each paragraph starts with an empty <w:r/>
#    that must be preserved).
Replace-ParagraphXml "Snake Wilds offers more potential for winning`r" `
    '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Adventure-themed slot with unique Snake Wilds feature</w:t></w:r></w:p>'

Replace-ParagraphXml "Detailed and quality graphics immerses players`r" `
    '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>High-quality graphics and immersive sound effects</w:t></w:r></w:p>'

Replace-ParagraphXml "Impressive sound effects enhancing overall experience`r" `
    '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Two types of free spins for increased winning potential</w:t></w:r></w:p>'

Replace-ParagraphXml "Free spin bonus round provides high payout potential`r" `
    '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Wide range of bet sizes and high RTP</w:t></w:r></w:p>'

# 3. "What we don't like" bullet list
Replace-ParagraphXml "High volatility level means less frequent payouts`r" `
    '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>High volatility may result in less frequent payouts</w:t></w:r></w:p>'

Replace-ParagraphXml "Limited number of paylines`r" `
    '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Similar theme to other popular slot games</w:t></w:r></w:p>'

# 4. Bold recap title near the end
Replace-ParagraphXml "Play Book of Souls II: El Dorado for Free - Review`r" `
    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Book of Souls II: El Dorado Free | Adventure-Themed Slot</w:t></w:r></w:p>'

# 5. Italic meta description paragraph (document's final paragraph)
Replace-ParagraphXml "Explore the ancient ruins of El Dorado in Book of Souls II. Enjoy unique features, such as Snake Wilds and two types of free spins for high payout potential. Play now for free.`r" `
    '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Book of Souls II: El Dorado and discover the adventure-themed slot with unique features. Play for free!</w:t></w:r></w:p>'
